$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Neue Messreihen mit Kondensator: M20-M23 (rows 37-40), D06-D19 shifted
# down by 7 rows (45-58, unchanged content) and new D20-D23 (59-62).
# Old rows 41-44 are vacated (their content now lives at rows 48-51).
$rowData = @(
    @(37, "M20", "2,7K", 100, 500, -5.5, "Schaltnetzteil", "40µF parallel"),
    @(38, "M21", "2,7K", 100, 500, -5.5, "Schaltnetzteil", "30µF parallel"),
    @(39, "M22", "2,7K", 100, 500, -5.5, "Schaltnetzteil", "20µF parallel"),
    @(40, "M23", "2,7K", 100, 500, -5.5, "Schaltnetzteil", "10µF parallel"),
    @(45, "D06", "18K", 570, 500, "ohne", "Schaltnetzteil", "Messung und Berechnung durch den µC, Übertragung der Leistungsdaten"),
    @(46, "D07", "18K", 570, 500, "ohne", "Schaltnetzteil", "Berechnete Leisuntgsdaten mit vorher gedrehtem Vorzeichen von i_temp"),
    @(47, "D08", "18K", 570, 500, "ohne", "Schaltnetzteil", "Kein Verstärker mit drehung VZ i_temp für positive P "),
    @(48, "D09", "18K", 570, 500, 4.7, "Schaltnetzteil", "nicht-invertierender Verstärker und -i_temp, Übersteuern von i!"),
    @(49, "D10", "18K", 570, 500, -4.7, "Schaltnetzteil", "invertierender Verstärker, Übersteuern von i bei Messung 12"),
    @(50, "D11", "18K", 570, 500, -4.7, 7805, "invertierender Verstärker"),
    @(51, "D12", "18K", 570, 500, 4.7, 7805, "nicht-invertierender Verstärker und -i_temp, Übersteuern von i! "),
    @(52, "D13", "18K", 570, 500, "ohne", 7805, "ohne Verstärker"),
    @(53, "D14", "18K", 570, 500, 3.2, 7805, "nicht-invertierender Verstärker"),
    @(54, "D15", "18K", 570, 500, -3.2, 7805, "invertierender Verstärker"),
    @(55, "D16", "18K", 570, 500, -3.2, "Schaltnetzteil", "invertierender Verstärker"),
    @(56, "D17", "18K", 570, 500, 3.2, "Schaltnetzteil", "nicht-invertierender Verstärker"),
    @(57, "D18", "18K", 570, 500, -5.5, "Schaltnetzteil", "3V LM317 Spannungs teiler von U nicht auf 3V ausgelegt"),
    @(58, "D19", "2,7K", 100, 500, -5.5, "Schaltnetzteil", "3V LM317 U jetzt wieder voll ausgesteuert"),
    @(59, "D20", "2,7K", 100, 500, -5.5, "Schaltnetzteil", "40µF parallel"),
    @(60, "D21", "2,7K", 100, 500, -5.5, "Schaltnetzteil", "30µF parallel"),
    @(61, "D22", "2,7K", 100, 500, -5.5, "Schaltnetzteil", "20µF parallel"),
    @(62, "D23", "2,7K", 100, 500, -5.5, "Schaltnetzteil", "10µF parallel")
)

foreach ($r in $rowData) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
    $ws.Range("D$rowNum").Value = $r[4]
    $ws.Range("E$rowNum").Value = $r[5]
    $ws.Range("F$rowNum").Value = $r[6]
    $ws.Range("G$rowNum").Value = $r[7]
}

# Rows 41:44 held the D09-D12 entries before the insert; that data now
# lives at rows 48-51, so the old slots are cleared out entirely.
$ws.Range("A41:G44").ClearContents()

# Restore the view: scrolled down to the new rows, cell A41 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A41").Select()
